# Data -> Data Tools -> Remove Duplicates
# The sheet has a header row in row 1 (B1:I1 populated, A1 empty) and data
# rows 2-48. Rows 46 and 47 are exact duplicates across all nine columns
# (A-I); removing the duplicate shifts row 48's data up into row 47 and
# leaves row 48 blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1:I48")
$rng.RemoveDuplicates(@(1,2,3,4,5,6,7,8,9), [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)

# The vacated last row keeps only the styled column G cell (column G carries
# an explicit column style), reset to the default style with no leftover
# value; columns H and I are fully cleared (no cell left behind at all).
$ws.Range("H48:I48").Clear()
$ws.Range("G48").Style = "Normal"
